# "Some changes in AAI"
# - Add a new column F ("unrelized") with boolean flag values for rows 2-26
# - Flip a handful of existing boolean values in columns A-E

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "unrelized"

# Corrected / toggled boolean values in the existing A:E grid
$ws.Range("B2").Value = $TRUE

$ws.Range("A3").Value = $TRUE
$ws.Range("D3").Value = $FALSE
$ws.Range("E3").Value = $FALSE

$ws.Range("B4").Value = $TRUE
$ws.Range("C4").Value = $TRUE

$ws.Range("A5").Value = $TRUE
$ws.Range("D5").Value = $FALSE
$ws.Range("E5").Value = $FALSE

$ws.Range("E13").Value = $TRUE

$ws.Range("E14").Value = $FALSE

$ws.Range("B15").Value = $TRUE
$ws.Range("C15").Value = $TRUE

$ws.Range("A16").Value = $TRUE
$ws.Range("D16").Value = $FALSE
$ws.Range("E16").Value = $FALSE

$ws.Range("B17").Value = $TRUE
$ws.Range("C17").Value = $TRUE

$ws.Range("B20").Value = $TRUE
$ws.Range("C20").Value = $TRUE

$ws.Range("A21").Value = $TRUE
$ws.Range("D21").Value = $FALSE
$ws.Range("E21").Value = $FALSE

$ws.Range("B22").Value = $TRUE
$ws.Range("C22").Value = $TRUE

$ws.Range("B25").Value = $TRUE
$ws.Range("C25").Value = $TRUE

$ws.Range("A26").Value = $TRUE
$ws.Range("D26").Value = $FALSE
$ws.Range("E26").Value = $FALSE

# New column F values, row by row (matching the final grid)
$fValues = @{
    2 = $TRUE
    3 = $FALSE
    4 = $TRUE
    5 = $FALSE
    6 = $FALSE
    7 = $FALSE
    8 = $FALSE
    9 = $FALSE
    10 = $FALSE
    11 = $FALSE
    12 = $FALSE
    13 = $FALSE
    14 = $FALSE
    15 = $TRUE
    16 = $FALSE
    17 = $TRUE
    18 = $FALSE
    19 = $FALSE
    20 = $TRUE
    21 = $FALSE
    22 = $TRUE
    23 = $FALSE
    24 = $FALSE
    25 = $TRUE
    26 = $FALSE
}

foreach ($r in $fValues.Keys) {
    $ws.Cells.Item($r, 6).Value = $fValues[$r]
}
